$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")
Write-Host $ws.Name
